$d = $word.ActiveDocument

# --- 1) Collapse the spell-check-split runs back into single runs. ---
# These five Find/Replace calls don't change the visible text; they just
# get rid of the <w:proofErr/> wrapped runs that fragmented the text,
# leaving each passage as a single run again.

$d.Content.Find.Execute("Existe la posición (Xd y Yd)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Existe la posición (Xd y Yd)", 2)

$d.Content.Find.Execute("Existe la posición (Xo y Yo)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Existe la posición (Xo y Yo)", 2)

$d.Content.Find.Execute("Altura del cañón (Ho), la separación de los cañones (d), la posición en Xo, la posición en Yo, el radio de destrucción del cañón 0.05*d, radio de neutralización 0.005*d", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Altura del cañón (Ho), la separación de los cañones (d), la posición en Xo, la posición en Yo, el radio de destrucción del cañón 0.05*d, radio de neutralización 0.005*d", 2)

$d.Content.Find.Execute("Altura del cañón (Hd), la separación de los cañones (d), la posición en Xd, la posición en Yd, el radio de destrucción del cañón 0.025*d", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Altura del cañón (Hd), la separación de los cañones (d), la posición en Xd, la posición en Yd, el radio de destrucción del cañón 0.025*d", 2)

$d.Content.Find.Execute("La clase cañón tendrá la posición en X, Y y la distancia de separación de los cañones y el rango de destrucción", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "La clase cañón tendrá la posición en X, Y y la distancia de separación de los cañones y el rango de destrucción", 2)

# --- 2) Turn the final empty paragraph into the new "NOTA" block. ---
# The document ends with two empty paragraphs; the last one becomes the
# bold "NOTA:" paragraph, and three more bold paragraphs are appended
# after it.

$last = $d.Paragraphs.Last
$last.Range.Text = "NOTA"
$last.Range.Font.Bold = $true
$last.Range.InsertAfter(":")

$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.Text = "Se podrán realizar ajustes al análisis en el momento del desarrollo del parcial, esto debido a que en los exámenes anteriores se a presentado la necesidad de realizar cambios en dicho análisis."
$last.Range.Font.Bold = $true

$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.Text = "Las variables van a tomarse de manera flotante"
$last.Range.Font.Bold = $true

$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.Text = "Para el desarrollo grafico se implementara"
$last.Range.Font.Bold = $true
$last.Range.InsertAfter("n diferentes clases que permitan el fácil desarrollo del mismo")
